$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A slightly (13.85546875 -> 14.85546875 in raw OOXML width units)
$ws.Columns.Item(1).ColumnWidth = 14

# Row 3
$ws.Range("A3").Value = 42600.792164351849
$ws.Range("B3").Value = "Bag"
$ws.Range("C3").Value = 6409
$ws.Range("D3").Value = 9671
$ws.Range("E3").Value = 1129
$ws.Range("F3").Value = 161
$ws.Range("G3").Value = 61
$ws.Range("H3").Value = 72
$ws.Range("I3").Value = 27
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 66
$ws.Range("M3").Value = 33

# Row 4
$ws.Range("A4").Value = 42600.794490740744
$ws.Range("B4").Value = "Bag"
$ws.Range("C4").Value = 8579
$ws.Range("D4").Value = 11114
$ws.Range("E4").Value = 1326
$ws.Range("F4").Value = 175
$ws.Range("G4").Value = 80
$ws.Range("H4").Value = 68
$ws.Range("I4").Value = 31
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 66
$ws.Range("M4").Value = 33

# Row 5
$ws.Range("A5").Value = 42600.830787037034
$ws.Range("B5").Value = "Bag"
$ws.Range("C5").Value = 6736
$ws.Range("D5").Value = 8115
$ws.Range("E5").Value = 956
$ws.Range("F5").Value = 111
$ws.Range("G5").Value = 60
$ws.Range("H5").Value = 64
$ws.Range("I5").Value = 34
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 66
$ws.Range("M5").Value = 33

# Row 6
$ws.Range("A6").Value = 42600.879189814812
$ws.Range("B6").Value = "Bag"
$ws.Range("C6").Value = 5830
$ws.Range("D6").Value = 6498
$ws.Range("E6").Value = 754
$ws.Range("F6").Value = 73
$ws.Range("G6").Value = 49
$ws.Range("H6").Value = 59
$ws.Range("I6").Value = 39
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 50
$ws.Range("M6").Value = 50

# Apply the date/time number format (style index 1 in the original workbook) to the new date cells
$ws.Range("A3:A6").NumberFormat = "m/d/yy h:mm"
